# TC01.xlsx - update login data set (positive and negative login)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New credentials replacing the previous sample data.
# Row 2: negative login (B2 set first so the "Admin@123" shared string
# takes the lowest available slot, matching password reuse in row 3)
$ws.Range("B2").Value = "Admin@123"
$ws.Range("A3").Value = "fsqa1tpn@gmail.com"
$ws.Range("A2").Value = "fsssacct3@gmail.com"
$ws.Range("B3").Value = "Admin@123"

# Move the active selection to D3, as left by the author when saving.
[void]$ws.Range("D3").Select()
